# Weekly update: insert a new week's price row for
# Agrícola del Norte S.A. de Arica - Jengibre, at the top of the data
# (row 3), pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 3:11 down to 4:12, preserving all their
# values/formatting, and opening up a blank row 3 for the new entry.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44616
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100114007
$ws.Range("G3").Value = "Jengibre"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 19500
$ws.Range("N3").Value = "$/caja 13 kilos"
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 1500
$ws.Range("Q3").Value = 13
$ws.Range("R3").Value = "Hortaliza"
